# ART_dict.xlsx edit: add INSTRUCTIONS_SINGLE and INSTRUCTIONS_SINGLE_PAGE
# rows, rename INSTRUCTIONS -> INSTRUCTIONS_PAIRS and PROMPT -> PROMPT_PAIRS,
# drop all the PRACTICE_ITEMS* / CONTINUE_MAIN_TEST / EXAMPLE* (example-page)
# rows that are no longer used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out all existing data rows (keep the header row 1) so we can lay the
# table back down in the new order without worrying about leftover cells.
$ws.Range("A2:C21").EntireRow.Delete()

# key / de / en triples, in final row order.
$data = @(
  @("TESTNAME", "Literatenquiz", "Artist Recognition Test"),
  @("INSTRUCTIONS_PAIRS", "Sie werden nun mehrere Paare von Namen sehen und müssen jeweils entscheiden, welcher der Personen eine Literat:in (Dichter:in, Schriftsteller:in, Dramatiker:in) ist.", "You will be presented with a set of name pairs for which you have to decide who is a literary writer (poet, novelist, playwright) and who is not."),
  @("INSTRUCTIONS_SINGLE", "Sie werden nun Namen präsentiert bekommen und müssen jeweils entscheiden, ob diese Persoen eine Literat:in (Dichter:in, Schriftsteller:in, Dramatiker:in) ist oder nicht.", "You will be presented a set of names for each you will have to decide if it belongs to a literary writer (poet, novelist, playwright) or not."),
  @("INSTRUCTIONS_SINGLE_PAGE", "Sie werden eine Liste mit Name sehen und sollen dort ankreuzen, welcher der Personen eine Literat:in (Dichter:in, Schriftsteller:in, Dramatiker:in) ist. Sie haben dazu {{time_out}} Sekunden Zeit.", "You will be presented with a list of names and you are asked to select all names which belong to a literary writer (poet, novelist, playwright)."),
  @("PROMPT_PAIRS", "Welcher der Personen ist eine Literat:in?<br/> Klicken Sie auf den Namen, sie haben {{time_out}} Sekunden Zeit zu antworten.", "Who of the persons is a literary writer? <br/> Click on the name, you have {{time_out}} seconds."),
  @("PROMPT_SINGLE", "Ist <b>{{name}}</b> eine Literat:in?<br/> Klicken Sie Ja oder Nein, sie haben {{time_out}} Sekunden Zeit zu antworten.", "Is <b>{{name}}</b> a literary writer? <br/> Click Yes or No, you have {{time_out}} seconds."),
  @("PROMPT_SINGLE_PAGE", "Bitte wählen Sie alle Literat:innen (Dichter:innen, Romanautor:innen, Dramatiker:innen)  aus der untenstehenden Liste aus.  Sie haben {{time_out}} Sekunden Zeit.", "Please select all literary writers (poets, novelists, playwrights). You have {{time_out}} seconds."),
  @("FEEDBACK", "Sie haben {{num_correct}} von {{num_items}} Fragen richtig beantwortet ({{perc_correct}}%).", "You answered {{num_correct}} out of {{num_items}} questions correctly ({{perc_correct}}%)."),
  @("FEEDBACK_SINGLE_PAGE", "Sie haben {{num_correct}} Literaten aus {{num_items}} Namen richtig erkannt ({{perc_correct}}%, Punkte: {{points}}).", "You answered {{num_correct}} out of {{num_items}} questions correctly ({{perc_correct}}%,  Points: {{points}})."),
  @("EXAMPLE_HEADER", "Beispiel {{page_no}} von {{num_pages}}", "Example {{page_no}} of {{num_pages}}"),
  @("PAGE_COUNTER", "Seite {{page_no}} von {{num_pages}}", "Page {{page_no}} of {{num_pages}}"),
  @("WELCOME", "Willkommen zum Literatenquiz!", "Welcome to the Artist Recognition Test!"),
  @("YES", "Ja", "Yes"),
  @("NO", "Nein", "No")
)

$row = 2
foreach ($item in $data) {
  $ws.Cells.Item($row, 1).Value = $item[0]
  $ws.Cells.Item($row, 2).Value = $item[1]
  $ws.Cells.Item($row, 3).Value = $item[2]
  $row = $row + 1
}

# Rows 3-5 (the long *_PAIRS / *_SINGLE / *_SINGLE_PAGE instruction rows) get
# a taller, wrapped "B" column, matching the other instruction-style rows.
foreach ($r in 3..5) {
  $ws.Range("B" + $r).WrapText = $true
  $ws.Rows.Item($r).RowHeight = 20.1
}

# Rows 11-12 (EXAMPLE_HEADER / PAGE_COUNTER) keep the top-vertical-align
# styling used for the other "meta" rows.
foreach ($r in 11..12) {
  $ws.Range("A" + $r + ":C" + $r).VerticalAlignment = -4160
}

$null = $ws.Range("C5").Select()
